$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$wsSteps = $wb.Worksheets.Item("Test Case Steps")

# --- Add new row 17 content (two new test cases combined into a single row, per diff) ---

# Column A (TCID) - style like A16 (fill + border)
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "TestCase_F16"

# Column B (Jira id) - style like C16 (border only, no fill)
$ws.Range("C16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "OPQA-231,OPQA-1100"

# Column C (Description) - style with border + wrap text, no fill (as used on sheet "Test Case Steps")
$wsSteps.Range("A2").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "Verify that Trending now section include articles and posts and able to navigate from tending now section and `nVerify that Maximum count on the trending list is 10"

# Column D (Runmode) - style like D16 (fill + border)
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Y"

# Column E (Results) - style like E16 (fill + border)
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "PASS"

$excel.CutCopyMode = $false

# Row height for the new row
$ws.Rows.Item(17).RowHeight = 30

# --- Column width adjustments: column B widens to fit new content ---
$ws.Columns.Item(2).AutoFit()

# --- Update selection to reflect new active cell ---
$ws.Activate()
$ws.Range("D14").Select()
